$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 11500
$ws.Range("I24").Value = 2000
$ws.Range("J24").Value = 21000
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 21000
$ws.Range("M24").Value = -1626
$ws.Range("N24").Value = -21748
$ws.Range("H32").Value = 1552.5714
$ws.Range("I32").Value = 1556.6154
$ws.Range("K32").Value = 1556.6154
$ws.Range("M32").Value = -1269.6154
$ws.Range("H63").Value = 6120.6665
$ws.Range("I63").Value = 4385
$ws.Range("K63").Value = 4385
$ws.Range("M63").Value = -3699
$ws.Range("H66").Value = 6120.6665
$ws.Range("I66").Value = 4385
$ws.Range("K66").Value = 21925
$ws.Range("M66").Value = -18493
$ws.Range("H100").Value = 11500
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 21000
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 21000
$ws.Range("M100").Value = -918
$ws.Range("N100").Value = -23164
$ws.Range("H122").Value = 5769.385
$ws.Range("I122").Value = 6427.143
$ws.Range("J122").Value = 5002
$ws.Range("K122").Value = 19281.429
$ws.Range("L122").Value = 15006
$ws.Range("M122").Value = -16831.429
$ws.Range("N122").Value = -19906

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -713
$ws.Range("N16").ClearContents()
$ws.Range("H50").Value = 22495
$ws.Range("J50").Value = 34990
$ws.Range("L50").Value = 34990
$ws.Range("N50").Value = -36240
$ws.Range("H51").Value = 50299
$ws.Range("J51").Value = 50299
$ws.Range("L51").Value = 50299
$ws.Range("N51").Value = -51771
$ws.Range("H60").Value = 34995
$ws.Range("J60").Value = 34995
$ws.Range("L60").Value = 34995
$ws.Range("N60").Value = -36017
$ws.Range("H61").Value = 50299
$ws.Range("J61").Value = 50299
$ws.Range("L61").Value = 50299
$ws.Range("N61").Value = -50995
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1170
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -19900
$ws.Range("H124").Value = 42900
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 42900
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 42900
$ws.Range("N124").Value = -52720
$ws.Range("H125").Value = 50000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 50000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("H127").Value = 100000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 100000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 4665.6665
$ws.Range("I132").Value = 3999
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 11997
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -9467
$ws.Range("N132").Value = -23057
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 60000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 60000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 40000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 40000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12550
$ws.Range("H123").Value = 70000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 70000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 70000
$ws.Range("N123").Value = -79800
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 1370.25
$ws.Range("I126").Value = 1370.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4110.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1640.75
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 60000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 60000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 60000
$ws.Range("N128").Value = -69960
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 99999.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 99999.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 99999.5
$ws.Range("N131").Value = -110079.5
$ws.Range("H132").Value = 3917.6667
$ws.Range("I132").Value = 3917.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11753.0001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9223.000100000001
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 1675.75
$ws.Range("I136").Value = 1675.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5027.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2477.25
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 50390
$ws.Range("I140").Value = 50390
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 50390
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -45210
$ws.Range("H141").Value = 40000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 40000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360
